# Auto-generated edit script: applies the Tonberry_Profits market-price refresh
# across all affected leve-profit sheets (columns H-N: market price / profit data).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 222.42857
$ws.Range("I55").Value = 231
$ws.Range("J55").Value = 201
$ws.Range("K55").Value = 231
$ws.Range("L55").Value = 201
$ws.Range("M55").Value = -17
$ws.Range("N55").Value = -629
$ws.Range("H62").Value = 1586.75
$ws.Range("I62").Value = 1527.7142
$ws.Range("K62").Value = 1527.7142
$ws.Range("M62").Value = -903.7141999999999
$ws.Range("H64").Value = 3167.5
$ws.Range("I64").Value = 2905.7144
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 2905.7144
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -2657.7144
$ws.Range("N64").Value = -5496
$ws.Range("H65").Value = 1586.75
$ws.Range("I65").Value = 1527.7142
$ws.Range("K65").Value = 7638.571
$ws.Range("M65").Value = -4518.571
$ws.Range("H67").Value = 3167.5
$ws.Range("I67").Value = 2905.7144
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 2905.7144
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -2047.7144
$ws.Range("N67").Value = -6716
$ws.Range("H92").Value = 821097.7
$ws.Range("I92").Value = 1026205.94
$ws.Range("K92").Value = 1026205.94
$ws.Range("M92").Value = -1024957.94
$ws.Range("H98").Value = 2003.3235
$ws.Range("I98").Value = 2035.138
$ws.Range("J98").Value = 1818.8
$ws.Range("K98").Value = 2035.138
$ws.Range("L98").Value = 1818.8
$ws.Range("M98").Value = -537.1379999999999
$ws.Range("N98").Value = -4814.8
$ws.Range("H116").Value = 8878.277
$ws.Range("I116").Value = 13161.556
$ws.Range("J116").Value = 4595
$ws.Range("K116").Value = 13161.556
$ws.Range("L116").Value = 4595
$ws.Range("M116").Value = -9719.556
$ws.Range("N116").Value = -11479
$ws.Range("H121").Value = 1159.25
$ws.Range("J121").Value = 1495.6666
$ws.Range("L121").Value = 4486.9998
$ws.Range("N121").Value = -7980.9998
$ws.Range("H122").Value = 2003.3235
$ws.Range("I122").Value = 2035.138
$ws.Range("J122").Value = 1818.8
$ws.Range("K122").Value = 6105.414
$ws.Range("L122").Value = 5456.4
$ws.Range("M122").Value = -3655.414
$ws.Range("N122").Value = -10356.4
$ws.Range("H132").Value = 7093676.5
$ws.Range("I132").Value = 8334728
$ws.Range("J132").Value = 1955
$ws.Range("K132").Value = 25004184
$ws.Range("L132").Value = 5865
$ws.Range("M132").Value = -25001654
$ws.Range("N132").Value = -10925
$ws.Range("H138").Value = 1822.2898
$ws.Range("J138").Value = 2648.7646
$ws.Range("L138").Value = 7946.293799999999
$ws.Range("N138").Value = -18226.2938
$ws.Range("H141").Value = 719602.2
$ws.Range("I141").Value = 824318.25
$ws.Range("K141").Value = 2472954.75
$ws.Range("M141").Value = -2467774.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 397553.28
$ws.Range("I2").Value = 428071.78
$ws.Range("J2").Value = 813
$ws.Range("K2").Value = 428071.78
$ws.Range("L2").Value = 813
$ws.Range("M2").Value = -427958.78
$ws.Range("N2").Value = -1039
$ws.Range("H32").Value = 3013.0115
$ws.Range("I32").Value = 2574.9023
$ws.Range("J32").Value = 10198
$ws.Range("K32").Value = 2574.9023
$ws.Range("L32").Value = 10198
$ws.Range("M32").Value = -2287.9023
$ws.Range("N32").Value = -10772
$ws.Range("H45").Value = 1766.85
$ws.Range("I45").Value = 1678.75
$ws.Range("K45").Value = 1678.75
$ws.Range("M45").Value = -1301.75
$ws.Range("H61").Value = 45457384
$ws.Range("I61").Value = 26317918
$ws.Range("J61").Value = 166674000
$ws.Range("K61").Value = 26317918
$ws.Range("L61").Value = 166674000
$ws.Range("M61").Value = -26317706
$ws.Range("N61").Value = -166674424
$ws.Range("H116").Value = 397553.28
$ws.Range("I116").Value = 428071.78
$ws.Range("J116").Value = 813
$ws.Range("K116").Value = 428071.78
$ws.Range("L116").Value = 813
$ws.Range("M116").Value = -425777.78
$ws.Range("N116").Value = -5401
$ws.Range("H132").Value = 1442.1552
$ws.Range("I132").Value = 1117.9048
$ws.Range("J132").Value = 2293.3125
$ws.Range("K132").Value = 3353.7144
$ws.Range("L132").Value = 6879.9375
$ws.Range("M132").Value = -823.7143999999998
$ws.Range("N132").Value = -11939.9375
$ws.Range("H136").Value = 45457384
$ws.Range("I136").Value = 26317918
$ws.Range("J136").Value = 166674000
$ws.Range("K136").Value = 78953754
$ws.Range("L136").Value = 500022000
$ws.Range("M136").Value = -78951204
$ws.Range("N136").Value = -500027100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 397553.28
$ws.Range("I3").Value = 428071.78
$ws.Range("J3").Value = 813
$ws.Range("K3").Value = 428071.78
$ws.Range("L3").Value = 813
$ws.Range("M3").Value = -427957.78
$ws.Range("N3").Value = -1041
$ws.Range("H20").Value = 1413.625
$ws.Range("I20").Value = 1172.5217
$ws.Range("J20").Value = 2029.7778
$ws.Range("K20").Value = 1172.5217
$ws.Range("L20").Value = 2029.7778
$ws.Range("M20").Value = -925.5217
$ws.Range("N20").Value = -2523.7778
$ws.Range("H105").Value = 2089
$ws.Range("I105").Value = 2052.4412
$ws.Range("K105").Value = 2052.4412
$ws.Range("M105").Value = -305.4412000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1784.8
$ws.Range("I99").Value = 1731
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1731
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -233
$ws.Range("N99").Value = -4996
$ws.Range("H126").Value = 1784.8
$ws.Range("I126").Value = 1731
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5193
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -2723
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 440.78946
$ws.Range("I5").Value = 311.57144
$ws.Range("K5").Value = 934.71432
$ws.Range("M5").Value = -822.71432
$ws.Range("H44").Value = 464.83334
$ws.Range("I44").Value = 445
$ws.Range("J44").Value = 474.75
$ws.Range("K44").Value = 1335
$ws.Range("L44").Value = 1424.25
$ws.Range("M44").Value = -937
$ws.Range("N44").Value = -2220.25
$ws.Range("H117").Value = 839.4286
$ws.Range("I117").Value = 625.6667
$ws.Range("J117").Value = 999.75
$ws.Range("K117").Value = 1877.0001
$ws.Range("L117").Value = 2999.25
$ws.Range("M117").Value = 1564.9999
$ws.Range("N117").Value = -9883.25
$ws.Range("H122").Value = 1094.3226
$ws.Range("J122").Value = 1298.2273
$ws.Range("L122").Value = 11684.0457
$ws.Range("N122").Value = -16584.0457
$ws.Range("H131").Value = 2051.2
$ws.Range("J131").Value = 2121.379
$ws.Range("L131").Value = 6364.137
$ws.Range("N131").Value = -16444.137
$ws.Range("H132").Value = 924.6667
$ws.Range("I132").Value = 812
$ws.Range("J132").Value = 1150
$ws.Range("K132").Value = 7308
$ws.Range("L132").Value = 10350
$ws.Range("M132").Value = -4778
$ws.Range("N132").Value = -15410
$ws.Range("H135").Value = 440.78946
$ws.Range("I135").Value = 311.57144
$ws.Range("K135").Value = 2804.14296
$ws.Range("M135").Value = -269.1429600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5000
$ws.Range("J80").Value = 5000
$ws.Range("L80").Value = 5000
$ws.Range("N80").Value = -6996
$ws.Range("H83").Value = 5000
$ws.Range("J83").Value = 5000
$ws.Range("L83").Value = 25000
$ws.Range("N83").Value = -34984
$ws.Range("H102").Value = 2468.818
$ws.Range("I102").Value = 2588
$ws.Range("K102").Value = 2588
$ws.Range("M102").Value = -966
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -99120

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 23169.334
$ws.Range("I25").Value = 2500
$ws.Range("J25").Value = 33504
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 33504
$ws.Range("M25").Value = -2270
$ws.Range("N25").Value = -33964
$ws.Range("H43").Value = 10676
$ws.Range("J43").Value = 10676
$ws.Range("L43").Value = 10676
$ws.Range("N43").Value = -11062
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").Value = ""
$ws.Range("H122").Value = 2698.7222
$ws.Range("J122").Value = 3924.875
$ws.Range("L122").Value = 11774.625
$ws.Range("N122").Value = -16674.625
$ws.Range("H132").Value = 1634.9818
$ws.Range("I132").Value = 1387.1
$ws.Range("J132").Value = 1932.44
$ws.Range("K132").Value = 4161.299999999999
$ws.Range("L132").Value = 5797.32
$ws.Range("M132").Value = -1631.299999999999
$ws.Range("N132").Value = -10857.32
$ws.Range("H136").Value = 2919.8948
$ws.Range("I136").Value = 1899.0385
$ws.Range("K136").Value = 5697.1155
$ws.Range("M136").Value = -3147.1155

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 12079522
$ws.Range("I136").Value = 14247180
$ws.Range("J136").Value = 2571.1428
$ws.Range("K136").Value = 42741540
$ws.Range("L136").Value = 7713.428400000001
$ws.Range("M136").Value = -42738990
$ws.Range("N136").Value = -12813.4284
